$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 526640.0600000001
$ws.Range("I2").Value = 833581.5
$ws.Range("K2").Value = 833581.5
$ws.Range("M2").Value = -833468.5
$ws.Range("H40").Value = 1505.6364
$ws.Range("J40").Value = 1128
$ws.Range("L40").Value = 1128
$ws.Range("N40").Value = -1478
$ws.Range("H51").Value = 2175
$ws.Range("H58").Value = 17858114
$ws.Range("I58").Value = 20833632
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 62500896
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -62500746
$ws.Range("N58").Value = -15300

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2877.8298
$ws.Range("I32").Value = 2600.558
$ws.Range("J32").Value = 5858.5
$ws.Range("K32").Value = 2600.558
$ws.Range("L32").Value = 5858.5
$ws.Range("M32").Value = -2313.558
$ws.Range("N32").Value = -6432.5
$ws.Range("H63").Value = 2903
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2903
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 2045
$ws.Range("I122").Value = 1751.1538
$ws.Range("K122").Value = 5253.4614
$ws.Range("M122").Value = -2803.4614

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 684.1667
$ws.Range("I64").Value = 494.625
$ws.Range("J64").Value = 1063.25
$ws.Range("K64").Value = 494.625
$ws.Range("L64").Value = 1063.25
$ws.Range("M64").Value = -269.625
$ws.Range("N64").Value = -1513.25
$ws.Range("H67").Value = 684.1667
$ws.Range("I67").Value = 494.625
$ws.Range("J67").Value = 1063.25
$ws.Range("K67").Value = 494.625
$ws.Range("L67").Value = 1063.25
$ws.Range("M67").Value = 285.375
$ws.Range("N67").Value = -2623.25
$ws.Range("H99").Value = 2886.3333
$ws.Range("I99").Value = 1220
$ws.Range("K99").Value = 1220
$ws.Range("M99").Value = 278
$ws.Range("H134").Value = 1924.3077
$ws.Range("I134").Value = 1641.7
$ws.Range("J134").Value = 2866.3333
$ws.Range("K134").Value = 4925.1
$ws.Range("L134").Value = 8598.999899999999
$ws.Range("M134").Value = -2390.1
$ws.Range("N134").Value = -13668.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1064.7142
$ws.Range("I22").Value = 514.1429000000001
$ws.Range("K22").Value = 514.1429000000001
$ws.Range("M22").Value = -164.1429000000001
$ws.Range("H132").Value = 2368.3333
$ws.Range("I132").Value = 2478.923
$ws.Range("J132").Value = 1649.5
$ws.Range("K132").Value = 7436.768999999999
$ws.Range("L132").Value = 4948.5
$ws.Range("M132").Value = -4906.768999999999
$ws.Range("N132").Value = -10008.5
$ws.Range("H134").Value = 2531.5
$ws.Range("I134").Value = 2559.65
$ws.Range("K134").Value = 7678.950000000001
$ws.Range("M134").Value = -5143.950000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1015
$ws.Range("J52").Value = 1015
$ws.Range("L52").Value = 3045
$ws.Range("N52").Value = -3577
$ws.Range("H70").Value = 530.6
$ws.Range("I70").Value = 530.6
$ws.Range("K70").Value = 1591.8
$ws.Range("M70").Value = -1276.8
$ws.Range("H73").Value = 530.6
$ws.Range("I73").Value = 530.6
$ws.Range("K73").Value = 1591.8
$ws.Range("M73").Value = -499.8000000000002
$ws.Range("H75").Value = 1609.875
$ws.Range("I75").Value = 479.66666
$ws.Range("J75").Value = 2288
$ws.Range("K75").Value = 1438.99998
$ws.Range("L75").Value = 6864
$ws.Range("M75").Value = -440.9999800000001
$ws.Range("N75").Value = -8860
$ws.Range("H78").Value = 1609.875
$ws.Range("I78").Value = 479.66666
$ws.Range("J78").Value = 2288
$ws.Range("K78").Value = 4316.99994
$ws.Range("L78").Value = 20592
$ws.Range("M78").Value = 675.0000600000003
$ws.Range("N78").Value = -30576
$ws.Range("H87").Value = 6205.4
$ws.Range("I87").Value = 5009
$ws.Range("J87").Value = 8000
$ws.Range("K87").Value = 15027
$ws.Range("L87").Value = 24000
$ws.Range("M87").Value = -13779
$ws.Range("N87").Value = -26496
$ws.Range("H90").Value = 6205.4
$ws.Range("I90").Value = 5009
$ws.Range("J90").Value = 8000
$ws.Range("K90").Value = 45081
$ws.Range("L90").Value = 72000
$ws.Range("M90").Value = -38841
$ws.Range("N90").Value = -84480
$ws.Range("H136").Value = 6875
$ws.Range("H140").Value = 2358.4285
$ws.Range("I140").Value = 2259.8333
$ws.Range("K140").Value = 6779.499899999999
$ws.Range("M140").Value = -1599.499899999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13090.546
$ws.Range("H73").Value = 13090.546
$ws.Range("H126").Value = 2584.7
$ws.Range("I126").Value = 2427.4443
$ws.Range("K126").Value = 7282.3329
$ws.Range("M126").Value = -4812.3329
$ws.Range("H132").Value = 2930.125
$ws.Range("I132").Value = 2240.6667
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 6722.000100000001
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -4192.000100000001
$ws.Range("N132").Value = -20055.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6104.4443
$ws.Range("I40").Value = 5867.5
$ws.Range("K40").Value = 5867.5
$ws.Range("M40").Value = -5731.5
$ws.Range("H61").Value = 3047.879
$ws.Range("I61").Value = 1938.8695
$ws.Range("J61").Value = 5598.6
$ws.Range("K61").Value = 1938.8695
$ws.Range("L61").Value = 5598.6
$ws.Range("M61").Value = -1736.8695
$ws.Range("N61").Value = -6002.6
$ws.Range("H82").Value = 1934.1538
$ws.Range("I82").Value = 811
$ws.Range("J82").Value = 4055.6667
$ws.Range("K82").Value = 811
$ws.Range("L82").Value = 4055.6667
$ws.Range("M82").Value = -450
$ws.Range("N82").Value = -4777.6667
$ws.Range("H85").Value = 1934.1538
$ws.Range("I85").Value = 811
$ws.Range("J85").Value = 4055.6667
$ws.Range("K85").Value = 811
$ws.Range("L85").Value = 4055.6667
$ws.Range("M85").Value = 437
$ws.Range("N85").Value = -6551.6667
$ws.Range("H93").Value = 3363.9119
$ws.Range("I93").Value = 763.5263
$ws.Range("J93").Value = 6657.7334
$ws.Range("K93").Value = 763.5263
$ws.Range("L93").Value = 6657.7334
$ws.Range("M93").Value = 484.4737
$ws.Range("N93").Value = -9153.733400000001
$ws.Range("H100").Value = 5535.92
$ws.Range("I100").Value = 4692.7856
$ws.Range("J100").Value = 6609
$ws.Range("K100").Value = 4692.7856
$ws.Range("L100").Value = 6609
$ws.Range("M100").Value = -4151.7856
$ws.Range("N100").Value = -7691
$ws.Range("H113").Value = 3047.879
$ws.Range("I113").Value = 1938.8695
$ws.Range("J113").Value = 5598.6
$ws.Range("K113").Value = 1938.8695
$ws.Range("L113").Value = 5598.6
$ws.Range("M113").Value = 231.1305
$ws.Range("N113").Value = -9938.6
$ws.Range("H122").Value = 7869.25
$ws.Range("I122").Value = 7869.25
$ws.Range("K122").Value = 23607.75
$ws.Range("M122").Value = -21157.75
$ws.Range("H136").Value = 2772.8215
$ws.Range("I136").Value = 1586.8125
$ws.Range("J136").Value = 4354.1665
$ws.Range("K136").Value = 4760.4375
$ws.Range("L136").Value = 13062.4995
$ws.Range("M136").Value = -2210.4375
$ws.Range("N136").Value = -18162.4995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13463.487
$ws.Range("J62").Value = 13815.527
$ws.Range("L62").Value = 13815.527
$ws.Range("N62").Value = -15063.527
$ws.Range("H64").Value = 142249.75
$ws.Range("J64").Value = 142249.75
$ws.Range("L64").Value = 142249.75
$ws.Range("N64").Value = -142745.75
$ws.Range("H65").Value = 13463.487
$ws.Range("J65").Value = 13815.527
$ws.Range("L65").Value = 69077.63499999999
$ws.Range("N65").Value = -75317.63499999999
$ws.Range("H67").Value = 142249.75
$ws.Range("J67").Value = 142249.75
$ws.Range("L67").Value = 142249.75
$ws.Range("N67").Value = -143965.75
$ws.Range("H122").Value = 4545.5454
$ws.Range("I122").Value = 2928.7144
$ws.Range("J122").Value = 7375
$ws.Range("K122").Value = 8786.143199999999
$ws.Range("L122").Value = 22125
$ws.Range("M122").Value = -6336.143199999999
$ws.Range("N122").Value = -27025
$ws.Range("H126").Value = 1042
$ws.Range("J126").Value = 913
$ws.Range("L126").Value = 2739
$ws.Range("N126").Value = -7679
$ws.Range("H132").Value = 4082.4167
$ws.Range("I132").Value = 3726.9
$ws.Range("J132").Value = 5860
$ws.Range("K132").Value = 11180.7
$ws.Range("L132").Value = 17580
$ws.Range("M132").Value = -8650.700000000001
$ws.Range("N132").Value = -22640
$ws.Range("H136").Value = 4271
$ws.Range("I136").Value = 4082.7
$ws.Range("K136").Value = 12248.1
$ws.Range("M136").Value = -9698.099999999999
